# Completed Working with the COUNT() Function
# Adds a new "COUNT" summary row (row 14) below MIN/MAX/AVERAGE, counting
# how many numeric entries are in each month's column (B:D), mirroring the
# existing MIN/MAX/AVERAGE rows' layout and fill pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row label
$ws.Range("A14").Value = "COUNT"

# B14 gets the typed-in formula; C14:D14 are then filled across from it,
# which is what produces the "shared formula" grouping Excel records for
# C14:D14 (matching the MIN/MAX/AVERAGE rows above).
$ws.Range("B14").Formula = "=COUNT(B4:B8)"
$ws.Range("C14:D14").Formula = "=COUNT(C4:C8)"

# Move the selection on, same as the recorded session did.
$ws.Range("I8").Select()
